$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 2034
    $ws.Range("F4").Value = 245
    $ws.Range("F6").Value = 6363
    $ws.Range("F7").Value = 237
}
